$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 231 (matching how the edit was made interactively) and delete it,
# shifting all rows below it up by one.
$ws.Rows.Item(231).Select()
$ws.Rows.Item(231).Delete()
